# Workbook/sheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each weekly block has:
#   - a "Hiếu / 08:00-16:00" row that must be cleared (kept blank)
#   - a "Long / 08:00-16:00 .. 08:00-17:00" row that is unchanged
#   - a brand-new row below, filled with "Hiếu / 13:00-19:00 / b / c / d / e / f / g"
$blocks = @(
    @{ Clear = 2;  New = 4 },
    @{ Clear = 8;  New = 10 },
    @{ Clear = 14; New = 16 },
    @{ Clear = 20; New = 22 },
    @{ Clear = 26; New = 28 },
    @{ Clear = 32; New = 34 }
)

$newRowValues = @("Hiếu", "13:00 - 19:00", "b", "c", "d", "e", "f", "g")

foreach ($block in $blocks) {
    $clearRow = $block.Clear
    $newRow = $block.New

    # Clear the old "Hiếu / 08:00-16:00" row (columns A-H), keeping its formatting
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($clearRow, $col).ClearContents()
    }

    # Fill in the new row below with Hiếu's updated shift and placeholder letters
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $newRowValues[$col - 1]
    }
}
